# Apply the 2022-05-23 data update to the "Fonds de solidarite" workbook.
# For each affected row, update nombre_aides (C), nombre_entreprises (D)
# and/or montant_total (E) to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 17;  C = 134742; D = $null;  E = 296794586 },
    @{ Row = 81;  C = 26163;  D = 2828;   E = 165085493 },
    @{ Row = 117; C = 19709;  D = $null;  E = 56504296 },
    @{ Row = 126; C = 5641;   D = $null;  E = 8171016 },
    @{ Row = 152; C = 126044; D = $null;  E = 715905660 },
    @{ Row = 168; C = 284979; D = $null;  E = 1209979525 },
    @{ Row = 170; C = 367345; D = $null;  E = 2845166254 },
    @{ Row = 171; C = 115135; D = $null;  E = 445869720 },
    @{ Row = 174; C = 357207; D = $null;  E = 1017494052 },
    @{ Row = 175; C = 125530; D = $null;  E = 812358049 },
    @{ Row = 179; C = 235690; D = 29336;  E = 812533024 },
    @{ Row = 180; C = 141469; D = $null;  E = 340196652 },
    @{ Row = 186; C = 21935;  D = $null;  E = 40057160 },
    @{ Row = 203; C = 13104;  D = $null;  E = 33016552 },
    @{ Row = 205; C = 11126;  D = $null;  E = 44116847 },
    @{ Row = 267; C = 84975;  D = $null;  E = 156519209 },
    @{ Row = 295; C = 91336;  D = 9957;   E = 552974072 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
